$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Permissions")

# --- Row 17: fill in the G/H columns for the already-present "Read transaction statuses" row ---
$ws.Range("G17").Value = "ReadTransactionStatus"
$ws.Range("H17").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,", ",$E$2,", ",$F$2,", ",$G$2,") values(N''",B17,"'', N''",C17,"'', ",IF(TRIM(D17)<>"","N''"&D17&"''","null"),", ",IF(TRIM(E17)<>"","N''"&E17&"''","null")," , getdate(), N''",G17,"'');")'

# --- New rows 18-21: notification permissions ---
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Read notifications"
$ws.Range("C18").Value = "Pregled notifikacija"
$ws.Range("G18").Value = "ReadNotification"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Edit notifications"
$ws.Range("C19").Value = "Promena notifikacija"
$ws.Range("G19").Value = "EditNotification"

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Insert notifications"
$ws.Range("C20").Value = "Dodavanje notifikacija"
$ws.Range("G20").Value = "InsertNotification"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Delete notifications"
$ws.Range("C21").Value = "Brisanje notifikacija"
$ws.Range("G21").Value = "DeleteNotification"

# Fill H18:H21 with the same shared formula pattern (extend the fill from H17)
$ws.Range("H17").AutoFill($ws.Range("H17:H21"))

# Update selection to match the recorded view state
$ws.Range("G25").Select()
